$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new session 13 row of data
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 30

# Update the selection to reflect the next empty row (as Excel would after data entry)
$ws.Range("B15").Select()
